$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Four new inventory rows for the TLC59711 breakout boards (Itead), appended
# right after the existing data (previously ending at row 27).
$newRows = @(
    @("Board", 11, "micro",  "micro (1)",  "N/A", "Itead", "N/A", "N/A", "N/A", "N/A", "Yes"),
    @("Board", 11, "accel",  "accel (2)",  "N/A", "Itead", "N/A", "N/A", "N/A", "N/A", "Yes"),
    @("Board", 12, "touch",  "touch (3)",  "N/A", "Itead", "N/A", "N/A", "N/A", "N/A", "Yes"),
    @("Board", 11, "power",  "power (4)",  "N/A", "Itead", "N/A", "N/A", "N/A", "N/A", "Yes")
)

$startRow = 28
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}

# Update the view state to match: scrolled down a couple rows and the
# selected cell moved to G25.
$excel.Goto($ws.Range("G25"))
